# reviewdb.xlsx: remove the duplicate "ronenchen27@gmail.com" review row (old row 2)
# and keep the "eligitel@gmail.com" review (old row 3), which becomes the sheet's
# only data row. Also fix up the hyperlinks so they point at the surviving row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a copy of the current (non-hyperlink) cell formatting used by the email /
# recovery columns so it can be restored later - deleting/recreating hyperlinks
# below causes Excel to stamp the cells with its builtin "Hyperlink" style.
$ws.Range("C2").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Drop the first data row (appid/keyword/ronenchen27.../danfogel100.../review#1).
# This shifts the remaining data row (the eligitel@gmail.com / ronenchen27@gmail.com
# review) up from row 3 into row 2.
$ws.Rows(2).Delete()

# The sheet's hyperlink list still points at the old (now stale) addresses
# C2/C3/D3, so rebuild it for the row that remains.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:eligitel@gmail.com", [Type]::Missing, [Type]::Missing, "eligitel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:ronenchen27@gmail.com", [Type]::Missing, [Type]::Missing, "ronenchen27@gmail.com")

# Restore the original look of those two cells (Hyperlinks.Add overwrote it).
$ws.Range("H1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("H1").Clear()

# Clean up the now-unused builtin Hyperlink cell style that Add() registered.
$wb.Styles.Item("Hyperlink").Delete()

# Match the saved selection/active cell of the edited workbook.
[void]$ws.Range("A2").Select()
